# Commit: "Como equilibrar 2 tarjetas"
#
# Change "y verticales respecto a la mesa" -> "y horizontales respecto a la
# mesa" inside step 1, splitting the single run that used to hold
# "y verticales respecto a la mesa" into three runs:
#   "y "  |  "horizontales"  |  " respecto a la mesa"
# The neighbouring runs ("entre si " before it, and the lone space after it)
# must stay exactly as they were.
#
# Plainly assigning to Range.Text on this engine re-normalises (merges) every
# run in the whole paragraph that happens to share the same effective
# character formatting, which would destroy the paragraph's existing run
# boundaries far beyond the word being edited. To avoid that, we briefly give
# the text we are about to touch (plus a one-character guard just past the
# end of the old run) a distinguishing character property (Bold) so the
# engine is forced to keep run boundaries exactly there; we edit the text
# while that marker is in place, then clear the temporary property again
# (set it back to "undefined" == wdUndefined == 9999999) so the saved
# run formatting ends up clean/unmarked, same as every other run around it.

$d = $word.ActiveDocument

$oldWord = "verticales"
$newWord = "horizontales"
$beforeWord = "y "
$oldPhrase = $beforeWord + $oldWord + " respecto a la mesa"

$wdUndefined = 9999999

$content = $d.Content
$phraseIdx = $content.Text.IndexOf($oldPhrase)
if ($phraseIdx -lt 0) {
    throw "Could not find phrase '$oldPhrase' in the document."
}
$phraseStart = $content.Start + $phraseIdx
$phraseEnd = $phraseStart + $oldPhrase.Length

# Guard character immediately after the original run. Marking it too keeps
# the post-edit re-normalisation from swallowing the runs that follow.
$d.Range($phraseEnd, $phraseEnd + 1).Bold = 1

# Locate & mark the word being replaced, then do the actual text swap.
$wordStart = $phraseStart + $beforeWord.Length
$wordEnd = $wordStart + $oldWord.Length
$target = $d.Range($wordStart, $wordEnd)
$target.Bold = 1
$target.Text = $newWord

# Clear the temporary Bold marks now that the run split has taken effect.
$newWordEnd = $wordStart + $newWord.Length
$d.Range($wordStart, $newWordEnd).Bold = $wdUndefined

$shift = $newWord.Length - $oldWord.Length
$guardStart = $phraseEnd + $shift
$d.Range($guardStart, $guardStart + 1).Bold = $wdUndefined
